# Apply the commit: "added one json for time bucket analysis"
#
# Net effect on this sheet: the data rows are re-ordered (a new/re-dated
# article bubbles to the top of the list), which rotates the (title, uri)
# pairs of the three existing data rows down by one:
#   old row2 (Statewide Blizzard / ohiohistory.org)      -> row3
#   old row3 (Pressure Records...   / weather.gov)       -> row4
#   old row4 (Blizzard of '78...    / miamistudent.net)  -> row2
#
# timestamp / historical distance / time bucket columns (B,C,D) are left
# untouched since every row already shares the same placeholder values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: title -------------------------------------------------
$ws.Range("A2").Value = "Blizzard of " + [char]0x2019 + "78 brought chilly attitudes to campus"
$ws.Range("A3").Value = "Statewide Blizzard"
$ws.Range("A4").Value = "Pressure Records: The October 26-27, 2010 Significant Extratropical Cyclone"

# --- Column E: uri (display text + hyperlink target) ------------------
$ws.Hyperlinks.Delete()

$ws.Range("E2").Value = "http://miamistudent.net/?p=110293"
$ws.Range("E3").Value = "https://web.archive.org/web/20060506165233/http://www.ohiohistory.org/etcetera/exhibits/swio/pages/content/1978_blizzard.htm"
$ws.Range("E4").Value = "http://www.weather.gov/dlh/101026_pressurerecords"

$ws.Hyperlinks.Add($ws.Range("E2"), "http://miamistudent.net/?p=110293") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://web.archive.org/web/20060506165233/http://www.ohiohistory.org/etcetera/exhibits/swio/pages/content/1978_blizzard.htm") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "http://www.weather.gov/dlh/101026_pressurerecords") | Out-Null
